$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "xJafy493"
$ws.Range("B2").Value = 23072455
$ws.Range("C2").Value = "neowgrv38"
$ws.Range("D2").Value = 't5#8Pn$Q'
$ws.Range("F2").Value = "nbmsDSUi"
$ws.Range("G2").Value = "YIOp"
